# The deck currently carries two theme parts:
#   theme1.xml -> "Office Theme" colours (only used by the Notes Master)
#   theme2.xml -> "Integral" colours     (used by the Slide Master / the
#                                          presentation's active design)
#
# The target edit swaps those two colour palettes: the deck's active
# design (the Slide Master's theme, i.e. theme2.xml) becomes the
# "Office Theme" palette, while the Notes-Master-only theme keeps the
# "Integral" palette (which now lives in theme1.xml).
#
# The PowerPoint object model edits themes through
# Slide.ThemeColorScheme (12 slots: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) which all slides/the slide master share, so setting
# it once repaints the presentation's single active theme part.
# (No RGB() helper is available here, so the values below are the
# plain decimal RGB-long equivalents of the target hex colours.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# New "Office Theme" palette
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
